$p = $ppt.ActivePresentation

function Set-BoldSuffix($textRange, [string]$needle) {
    $t = $textRange.Text
    $idx = $t.IndexOf($needle)
    if ($idx -ge 0) {
        $sub = $textRange.Characters($idx + 1, $needle.Length)
        $sub.Font.Bold = $true
    }
}

# --- Slide 11: Title - "...SR Paths" -> "...SR " + bold "LSPs" ---
$s11 = $p.Slides.Item(11)
$titleTr = $s11.Shapes.Item(1).TextFrame.TextRange
$full = $titleTr.Text
$idx = $full.IndexOf("Paths")
if ($idx -ge 0) {
    $old = $titleTr.Characters($idx + 1, 5)
    $old.Text = "LSPs"
}
Set-BoldSuffix $s11.Shapes.Item(1).TextFrame.TextRange "LSPs"

# --- Slide 13: Requirements/Scope bullets - bold trailing "SR LSPs" ---
$s13 = $p.Slides.Item(13)
$contentTr = $s13.Shapes.Item(2).TextFrame.TextRange

$para2 = $contentTr.Paragraphs(2, 1)
Set-BoldSuffix $para2 "SR LSPs"

$para3 = $contentTr.Paragraphs(3, 1)
Set-BoldSuffix $para3 "SR LSPs"

$para5 = $contentTr.Paragraphs(5, 1)
Set-BoldSuffix $para5 "SR LSPs"

# --- Slide 14: rename Tunnel-AD / Tunnel-DA boxes ---
$s14 = $p.Slides.Item(14)
for ($i = 1; $i -le $s14.Shapes.Count; $i++) {
    $shp = $s14.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $txt = $shp.TextFrame.TextRange.Text
        if ($txt -eq "Tunnel-AD") {
            $shp.TextFrame.TextRange.Text = "SR-Policy-AD"
        } elseif ($txt -eq "Tunnel-DA") {
            $shp.TextFrame.TextRange.Text = "SR-Policy-DA"
        }
    }
}
